# Update the Omg-Tnfrsf1b LR-pairs sheet to reflect the re-run of the NATMI
# pipeline per Dr Hou's advice: ligand/receptor "expressing cells" counts
# (columns E and K) go from 1 to 3 for every data row, with every downstream
# derived expression/specificity/edge-weight column recomputed accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns, in order, that receive new values for each data row (2..17).
$columns = @("E", "G", "H", "I", "J", "K", "M", "N", "O", "P", "Q", "R", "S", "T")

# Each inner array is: row number, followed by the new value for each of the
# columns above (same order).
$updates = @(
    @(2, 3, 1.840954666666667, 5.522864, 0.1665577163991766, 0.1665577163991766, 3, 11.84297166666667, 35.528915, 0.1201507426883529, 0.1201507426883529, 21.80237395695111, 196.22136561256, 0.02001203332583712, 0.02001203332583713),
    @(3, 3, 1.840954666666667, 5.522864, 0.1665577163991766, 0.1665577163991766, 3, 10.940628, 32.821884, 0.1109961770302011, 0.1109961770302011, 20.141200172864, 181.270801555776, 0.01848726977518903, 0.01848726977518904),
    @(4, 3, 1.840954666666667, 5.522864, 0.1665577163991766, 0.1665577163991766, 3, 73.27090866666667, 219.812726, 0.7433568483938062, 0.7433568483938062, 134.8884212408071, 1213.995791167264, 0.1238118191381613, 0.1238118191381613),
    @(5, 3, 1.840954666666667, 5.522864, 0.1665577163991766, 0.1665577163991766, 3, 2.513102666666667, 7.539308, 0.02549623188763971, 0.02549623188763971, 4.626508082012445, 41.638572738112, 0.004246594159989137, 0.004246594159989138),
    @(6, 3, 3.899801666666667, 11.699405, 0.3528289271705964, 0.3528289271705964, 3, 11.84297166666667, 35.528915, 0.1201507426883529, 0.1201507426883529, 46.18524064395277, 415.667165795575, 0.04239265764148194, 0.04239265764148194),
    @(7, 3, 3.899801666666667, 11.699405, 0.3528289271705964, 0.3528289271705964, 3, 10.940628, 32.821884, 0.1109961770302011, 0.1109961770302011, 42.66627930878, 383.99651377902, 0.03916266206160345, 0.03916266206160345),
    @(8, 3, 3.899801666666667, 11.699405, 0.3528289271705964, 0.3528289271705964, 3, 73.27090866666667, 219.812726, 0.7433568483938062, 0.7433568483938062, 285.7420117364478, 2571.67810562803, 0.2622777993237023, 0.2622777993237023),
    @(9, 3, 3.899801666666667, 11.699405, 0.3528289271705964, 0.3528289271705964, 3, 2.513102666666667, 7.539308, 0.02549623188763971, 0.02549623188763971, 9.800601967971112, 88.20541771174001, 0.008995808143808669, 0.008995808143808669),
    @(10, 3, 2.952690333333333, 8.858071000000001, 0.2671403962621152, 0.2671403962621152, 3, 11.84297166666667, 35.528915, 0.1201507426883529, 0.1201507426883529, 34.96862795810722, 314.717651622965, 0.03209711701295404, 0.03209711701295405),
    @(11, 3, 2.952690333333333, 8.858071000000001, 0.2671403962621152, 0.2671403962621152, 3, 10.940628, 32.821884, 0.1109961770302011, 0.1109961770302011, 32.304286536196, 290.738578825764, 0.02965156271542781, 0.02965156271542782),
    @(12, 3, 2.952690333333333, 8.858071000000001, 0.2671403962621152, 0.2671403962621152, 3, 73.27090866666667, 219.812726, 0.7433568483938062, 0.7433568483938062, 216.3463037346162, 1947.116733611546, 0.1985806430440785, 0.1985806430440785),
    @(13, 3, 2.952690333333333, 8.858071000000001, 0.2671403962621152, 0.2671403962621152, 3, 2.513102666666667, 7.539308, 0.02549623188763971, 0.02549623188763971, 7.420413950540889, 66.78372555486801, 0.006811073489654849, 0.00681107348965485),
    @(14, 3, 2.359506666666667, 7.07852, 0.213472960168112, 0.213472960168112, 3, 11.84297166666667, 35.528915, 0.1201507426883529, 0.1201507426883529, 27.94357060064444, 251.4921354058, 0.02564893470807983, 0.02564893470807984),
    @(15, 3, 2.359506666666667, 7.07852, 0.213472960168112, 0.213472960168112, 3, 10.940628, 32.821884, 0.1109961770302011, 0.1109961770302011, 25.81448470352, 232.33036233168, 0.02369468247798082, 0.02369468247798083),
    @(16, 3, 2.359506666666667, 7.07852, 0.213472960168112, 0.213472960168112, 3, 73.27090866666667, 219.812726, 0.7433568483938062, 0.7433568483938062, 172.8831974717245, 1555.94877724552, 0.1586865868878642, 0.1586865868878642),
    @(17, 3, 2.359506666666667, 7.07852, 0.213472960168112, 0.213472960168112, 3, 2.513102666666667, 7.539308, 0.02549623188763971, 0.02549623188763971, 5.929682496017779, 53.36714246416, 0.005442756094187057, 0.005442756094187057)
)

foreach ($rowData in $updates) {
    $rowNum = $rowData[0]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $colLetter = $columns[$i]
        $value = $rowData[$i + 1]
        $ws.Range("$colLetter$rowNum").Value = $value
    }
}
